# Fix source documents that used Heading 1 for normal text.
#
# The "Prerequisite tasks" section of this document has a body paragraph
# ("None") that was incorrectly styled as Heading 1 (with a stray
# bookmark and leftover direct-formatting / an empty trailing run).
# Turn it into a plain paragraph: no heading style, no bookmark, a
# single run with the text "None" and only the rtl run property.

$d = $word.ActiveDocument

# Locate the paragraph: the one whose text is exactly "None" (styled as
# Heading 1 in the source doc -- this is the bug being fixed).
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $text = $p.Range.Text
    if ($text -eq "None`r" -and $p.Style.NameLocal -eq "Heading 1") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $range = $target.Range
    $xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:pPr><w:rPr/></w:pPr>' +
           '<w:r><w:rPr><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">None</w:t></w:r>' +
           '</w:p>'
    $range.InsertXML($xml)
}
